# change stock price data to go back 5 years
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = 0.005
$ws.Range("B3").Value  = 0.0500000000000004
$ws.Range("B4").Value  = 0.0499999999999998
$ws.Range("B5").Value  = 0.00500000000000005
$ws.Range("B6").Value  = 0.0150000000000016
$ws.Range("B7").Value  = 0.005
$ws.Range("B8").Value  = 0.00499999999999982
$ws.Range("B9").Value  = 0.00500000000000028
$ws.Range("B10").Value = 0.0499999999999994
$ws.Range("B11").Value = 0.00500000000000009
$ws.Range("B12").Value = 0.00500000000000002
$ws.Range("B13").Value = 0.05
$ws.Range("B14").Value = 0.0500000000000004
$ws.Range("B15").Value = 0.0499999999999995
$ws.Range("B16").Value = 0.005
$ws.Range("B17").Value = 0.00500000000000003
$ws.Range("B19").Value = 0.005
$ws.Range("B20").Value = 0.00500000000000001
$ws.Range("B21").Value = 0.0500000000000003
$ws.Range("B22").Value = 0.0500000000000001
$ws.Range("B23").Value = 0.05
$ws.Range("B24").Value = 0.00500000000000003
$ws.Range("B26").Value = 0.0500000000000004
$ws.Range("B27").Value = 0.005
$ws.Range("B28").Value = 0.00500000000000014
$ws.Range("B29").Value = 0.00500000000000002
$ws.Range("B30").Value = 0.149999999999998
$ws.Range("B31").Value = 0.0499999999999998
$ws.Range("B32").Value = 0.2
